$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8164013333333333
$ws.Range("H2").Value = 2.449204
$ws.Range("I2").Value = 0.05618115571687973
$ws.Range("J2").Value = 0.05618115571687973
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 62.93598963536444
$ws.Range("R2").Value = 566.42390671828
$ws.Range("S2").Value = 0.01350491988614786
$ws.Range("T2").Value = 0.01350491988614786
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8164013333333333
$ws.Range("H3").Value = 2.449204
$ws.Range("I3").Value = 0.05618115571687973
$ws.Range("J3").Value = 0.05618115571687973
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 82.93007791898312
$ws.Range("R3").Value = 746.370701270848
$ws.Range("S3").Value = 0.01779528795744154
$ws.Range("T3").Value = 0.01779528795744154
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8164013333333333
$ws.Range("H4").Value = 2.449204
$ws.Range("I4").Value = 0.05618115571687973
$ws.Range("J4").Value = 0.05618115571687973
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 115.9508601807858
$ws.Range("R4").Value = 1043.557741627072
$ws.Range("S4").Value = 0.02488094787329034
$ws.Range("T4").Value = 0.02488094787329034
$ws.Range("I5").Value = 0.8862323361798529
$ws.Range("J5").Value = 0.8862323361798529
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 992.7867878941134
$ws.Range("R5").Value = 8935.08109104702
$ws.Range("S5").Value = 0.2130340066505007
$ws.Range("T5").Value = 0.2130340066505007
$ws.Range("I6").Value = 0.8862323361798529
$ws.Range("J6").Value = 0.8862323361798529
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2807126236240503
$ws.Range("T6").Value = 0.2807126236240502
$ws.Range("I7").Value = 0.8862323361798529
$ws.Range("J7").Value = 0.8862323361798529
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.392485705905302
$ws.Range("T7").Value = 0.392485705905302
$ws.Range("G8").Value = 0.8368233333333334
$ws.Range("I8").Value = 0.05758650810326746
$ws.Range("J8").Value = 0.05758650810326746
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 64.51031188087778
$ws.Range("R8").Value = 580.5928069279
$ws.Range("S8").Value = 0.01384274083603392
$ws.Range("T8").Value = 0.01384274083603392
$ws.Range("G9").Value = 0.8368233333333334
$ws.Range("I9").Value = 0.05758650810326746
$ws.Range("J9").Value = 0.05758650810326746
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("R9").Value = 765.0409089726401
$ws.Range("S9").Value = 0.01824043099656797
$ws.Range("T9").Value = 0.01824043099656797
$ws.Range("G10").Value = 0.8368233333333334
$ws.Range("I10").Value = 0.05758650810326746
$ws.Range("J10").Value = 0.05758650810326746
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.02550333627066557
$ws.Range("T10").Value = 0.02550333627066557
